# Scheduled market-data refresh: update crafted-item price/profit figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns) across
# the ALC, ARM, CRP, CUL and LTW sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 17
$ws.Range("H17").Value = 1448.196
$ws.Range("J17").Value = 1448.196
$ws.Range("L17").Value = 4344.588
$ws.Range("N17").Value = -4680.588

# Row 18
$ws.Range("H18").Value = 842.1539
$ws.Range("I18").Value = 842.1539
$ws.Range("K18").Value = 842.1539
$ws.Range("M18").Value = -558.1539

# Row 33
$ws.Range("H33").Value = 29412242
$ws.Range("I33").Value = 201.14285
$ws.Range("K33").Value = 201.14285
$ws.Range("M33").Value = 27.85714999999999

# Row 48
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584
$ws.Range("M48").ClearContents()

# Row 56
$ws.Range("H56").Value = 3000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 3000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068
$ws.Range("M56").ClearContents()

# Row 64
$ws.Range("H64").Value = 3215.739
$ws.Range("I64").Value = 2772.4
$ws.Range("J64").Value = 3338.889
$ws.Range("K64").Value = 2772.4
$ws.Range("L64").Value = 3338.889
$ws.Range("M64").Value = -2524.4
$ws.Range("N64").Value = -3834.889

# Row 67
$ws.Range("H67").Value = 3215.739
$ws.Range("I67").Value = 2772.4
$ws.Range("J67").Value = 3338.889
$ws.Range("K67").Value = 2772.4
$ws.Range("L67").Value = 3338.889
$ws.Range("M67").Value = -1914.4
$ws.Range("N67").Value = -5054.889

# Row 76
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 3100
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 3100
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -3730

# Row 79
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 3100
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 3100
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -5284

# Row 111
$ws.Range("H111").Value = 3060.6316
$ws.Range("I111").Value = 3102.9
$ws.Range("J111").Value = 3013.6667
$ws.Range("K111").Value = 9308.700000000001
$ws.Range("L111").Value = 9041.000100000001
$ws.Range("M111").Value = -6241.700000000001
$ws.Range("N111").Value = -15175.0001

# Row 137
$ws.Range("H137").Value = 1450.8276
$ws.Range("I137").Value = 1171.4166
$ws.Range("K137").Value = 3514.2498
$ws.Range("M137").Value = -964.2498000000001

$ws = $wb.Worksheets("ARM")
# Row 61
$ws.Range("H61").Value = 1372.5946
$ws.Range("I61").Value = 1314.6666
$ws.Range("J61").Value = 1850.5
$ws.Range("K61").Value = 1314.6666
$ws.Range("L61").Value = 1850.5
$ws.Range("M61").Value = -1102.6666
$ws.Range("N61").Value = -2274.5

# Row 74
$ws.Range("H74").Value = 32460.086
$ws.Range("I74").Value = 46395.184
$ws.Range("J74").Value = 8877.615
$ws.Range("K74").Value = 46395.184
$ws.Range("L74").Value = 8877.615
$ws.Range("M74").Value = -45521.184
$ws.Range("N74").Value = -10625.615

# Row 77
$ws.Range("H77").Value = 32460.086
$ws.Range("I77").Value = 46395.184
$ws.Range("J77").Value = 8877.615
$ws.Range("K77").Value = 231975.92
$ws.Range("L77").Value = 44388.075
$ws.Range("M77").Value = -227607.92
$ws.Range("N77").Value = -53124.075

# Row 88
$ws.Range("H88").Value = 2415.5
$ws.Range("J88").Value = 2456
$ws.Range("L88").Value = 2456
$ws.Range("N88").Value = -3268

# Row 91
$ws.Range("H91").Value = 2415.5
$ws.Range("J91").Value = 2456
$ws.Range("L91").Value = 2456
$ws.Range("N91").Value = -5264

# Row 136
$ws.Range("H136").Value = 1372.5946
$ws.Range("I136").Value = 1314.6666
$ws.Range("J136").Value = 1850.5
$ws.Range("K136").Value = 3943.9998
$ws.Range("L136").Value = 5551.5
$ws.Range("M136").Value = -1393.9998
$ws.Range("N136").Value = -10651.5

$ws = $wb.Worksheets("CRP")
# Row 31
$ws.Range("H31").Value = 23811324
$ws.Range("I31").Value = 40001240
$ws.Range("J31").Value = 2623.7646
$ws.Range("K31").Value = 40001240
$ws.Range("L31").Value = 2623.7646
$ws.Range("M31").Value = -40000945
$ws.Range("N31").Value = -3213.7646

# Row 34
$ws.Range("H34").Value = 23811324
$ws.Range("I34").Value = 40001240
$ws.Range("J34").Value = 2623.7646
$ws.Range("K34").Value = 40001240
$ws.Range("L34").Value = 2623.7646
$ws.Range("M34").Value = -40001038
$ws.Range("N34").Value = -3027.7646

# Row 60
$ws.Range("H60").Value = 16998
$ws.Range("J60").Value = 16998
$ws.Range("L60").Value = 16998
$ws.Range("N60").Value = -18020

# Row 86
$ws.Range("H86").Value = 5821.6577
$ws.Range("I86").Value = 4454.35
$ws.Range("J86").Value = 7340.8887
$ws.Range("K86").Value = 4454.35
$ws.Range("L86").Value = 7340.8887
$ws.Range("M86").Value = -3331.35
$ws.Range("N86").Value = -9586.8887

# Row 89
$ws.Range("H89").Value = 5821.6577
$ws.Range("I89").Value = 4454.35
$ws.Range("J89").Value = 7340.8887
$ws.Range("K89").Value = 22271.75
$ws.Range("L89").Value = 36704.4435
$ws.Range("M89").Value = -16655.75
$ws.Range("N89").Value = -47936.4435

$ws = $wb.Worksheets("CUL")
# Row 39
$ws.Range("H39").Value = 1868.1818
$ws.Range("I39").Value = 550
$ws.Range("J39").Value = 2966.6667
$ws.Range("K39").Value = 1650
$ws.Range("L39").Value = 8900.000100000001
$ws.Range("M39").Value = -1356
$ws.Range("N39").Value = -9488.000100000001

# Row 51
$ws.Range("H51").Value = 2870
$ws.Range("I51").Value = 920
$ws.Range("J51").Value = 3520
$ws.Range("K51").Value = 2760
$ws.Range("L51").Value = 10560
$ws.Range("M51").Value = -2300
$ws.Range("N51").Value = -11480

# Row 131
$ws.Range("H131").Value = 891.36
$ws.Range("I131").Value = 611.25
$ws.Range("J131").Value = 915.7174
$ws.Range("K131").Value = 1833.75
$ws.Range("L131").Value = 2747.1522
$ws.Range("M131").Value = 3206.25
$ws.Range("N131").Value = -12827.1522

$ws = $wb.Worksheets("LTW")
# Row 136
$ws.Range("H136").Value = 1442.68
$ws.Range("I136").Value = 1253.85
$ws.Range("J136").Value = 2198
$ws.Range("K136").Value = 3761.55
$ws.Range("L136").Value = 6594
$ws.Range("M136").Value = -1211.55
$ws.Range("N136").Value = -11694
